$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, C, E are set directly since they never look like numbers.
# Column D ("Price") often contains values that look numeric (e.g. "1.00",
# "0.658"); Excel auto-infers those as numbers on a plain .Value assignment,
# which would silently change the cell type/representation. To preserve the
# original text-cell semantics we force text formatting, assign the value,
# then restore the default "Normal" style so no stray formatting is left
# behind on the cell.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '66.506.52'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +6.68%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.554.72'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +3.36%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '418.44'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.09%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '130.72'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.72%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.658'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +5.34%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.546.00'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +3.25%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.782'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +7.91%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.176'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +25.56%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0000309'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +42.04%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '43.24'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.24%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '10.04'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +6.22%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.124.30'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +3.70%  '
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '20.46'
$ws.Range('D17').Style = "Normal"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.569.67'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +3.49%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.12'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +4.80%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.56'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.94%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '66.552.43'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +6.71%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '451.12'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -5.39%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '90.35'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.93%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.22'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -2.00%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '13.14'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.27%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.38'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.31%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.98'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -6.13%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '34.44'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +3.57%  '
$ws.Range('E29').Value = '  +0.82%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '12.46'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +4.25%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.78'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +4.62%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.118'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +5.47%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.34'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -3.65%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.161'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -3.32%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.998'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.27%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '39.03'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -4.23%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '57.10'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.40%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0₃0801'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +45.75%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0501'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.40%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.148'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +10.32%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.78'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +3.57%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.02'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '149.09'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +2.71%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '4.41'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.95%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.26'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.42%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.310'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.63%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.00'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.53%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.33'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -3.91%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.145'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +3.58%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '15.55'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -5.24%  '
